$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 21090
$ws.Range("J57").Value = 21090
$ws.Range("L57").Value = 63270
$ws.Range("N57").Value = -64268
$ws.Range("H129").Value = 950669.9
$ws.Range("J129").Value = 1425872.5
$ws.Range("L129").Value = 4277617.5
$ws.Range("N129").Value = -4287617.5
$ws.Range("H132").Value = 1176.5
$ws.Range("I132").Value = 1239.6888
$ws.Range("J132").Value = 770.2857
$ws.Range("K132").Value = 3719.0664
$ws.Range("L132").Value = 2310.8571
$ws.Range("M132").Value = -1189.0664
$ws.Range("N132").Value = -7370.8571
$ws.Range("H137").Value = 889.5625
$ws.Range("I137").Value = 855.0833
$ws.Range("J137").Value = 993
$ws.Range("K137").Value = 2565.2499
$ws.Range("L137").Value = 2979
$ws.Range("M137").Value = -15.2498999999998
$ws.Range("N137").Value = -8079
$ws.Range("H138").Value = 1327.86
$ws.Range("I138").Value = 613.9167
$ws.Range("J138").Value = 3163.7144
$ws.Range("K138").Value = 1841.7501
$ws.Range("L138").Value = 9491.143199999999
$ws.Range("M138").Value = 3298.2499
$ws.Range("N138").Value = -19771.1432
$ws.Range("H141").Value = 3339.4792
$ws.Range("I141").Value = 793
$ws.Range("J141").Value = 7583.6113
$ws.Range("K141").Value = 2379
$ws.Range("L141").Value = 22750.8339
$ws.Range("M141").Value = 2801
$ws.Range("N141").Value = -33110.8339

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 5006.5
$ws.Range("J15").Value = 5006.5
$ws.Range("L15").Value = 5006.5
$ws.Range("N15").Value = -5706.5
$ws.Range("H32").Value = 3698.35
$ws.Range("I32").Value = 2862.8481
$ws.Range("J32").Value = 6841.4287
$ws.Range("K32").Value = 2862.8481
$ws.Range("L32").Value = 6841.4287
$ws.Range("M32").Value = -2575.8481
$ws.Range("N32").Value = -7415.4287
$ws.Range("H45").Value = 2070
$ws.Range("I45").Value = 1337.5
$ws.Range("K45").Value = 1337.5
$ws.Range("M45").Value = -960.5
$ws.Range("H61").Value = 1295.5172
$ws.Range("I61").Value = 875.3333
$ws.Range("J61").Value = 1983.091
$ws.Range("K61").Value = 875.3333
$ws.Range("L61").Value = 1983.091
$ws.Range("M61").Value = -663.3333
$ws.Range("N61").Value = -2407.091
$ws.Range("H63").Value = 1777.3158
$ws.Range("I63").Value = 1697.8235
$ws.Range("J63").Value = 2453
$ws.Range("K63").Value = 1697.8235
$ws.Range("L63").Value = 2453
$ws.Range("M63").Value = -1011.8235
$ws.Range("N63").Value = -3825
$ws.Range("H66").Value = 1777.3158
$ws.Range("I66").Value = 1697.8235
$ws.Range("J66").Value = 2453
$ws.Range("K66").Value = 8489.1175
$ws.Range("L66").Value = 12265
$ws.Range("M66").Value = -5057.1175
$ws.Range("N66").Value = -19129
$ws.Range("H74").Value = 933.6739
$ws.Range("I74").Value = 850.4545000000001
$ws.Range("J74").Value = 1144.9231
$ws.Range("K74").Value = 850.4545000000001
$ws.Range("L74").Value = 1144.9231
$ws.Range("M74").Value = 23.54549999999995
$ws.Range("N74").Value = -2892.9231
$ws.Range("H77").Value = 933.6739
$ws.Range("I77").Value = 850.4545000000001
$ws.Range("J77").Value = 1144.9231
$ws.Range("K77").Value = 4252.2725
$ws.Range("L77").Value = 5724.6155
$ws.Range("M77").Value = 115.7275
$ws.Range("N77").Value = -14460.6155
$ws.Range("H122").Value = 876.5714
$ws.Range("I122").Value = 784
$ws.Range("K122").Value = 2352
$ws.Range("M122").Value = 98
$ws.Range("H132").Value = 1445.7561
$ws.Range("I132").Value = 1186
$ws.Range("J132").Value = 1718.5
$ws.Range("K132").Value = 3558
$ws.Range("L132").Value = 5155.5
$ws.Range("M132").Value = -1028
$ws.Range("N132").Value = -10215.5
$ws.Range("H136").Value = 1295.5172
$ws.Range("I136").Value = 875.3333
$ws.Range("J136").Value = 1983.091
$ws.Range("K136").Value = 2625.9999
$ws.Range("L136").Value = 5949.272999999999
$ws.Range("M136").Value = -75.9998999999998
$ws.Range("N136").Value = -11049.273

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 7699.2856
$ws.Range("I8").Value = 1395
$ws.Range("J8").Value = 8750
$ws.Range("K8").Value = 1395
$ws.Range("L8").Value = 8750
$ws.Range("M8").Value = -1255
$ws.Range("N8").Value = -9030
$ws.Range("H70").Value = 100459
$ws.Range("J70").Value = 100459
$ws.Range("L70").Value = 100459
$ws.Range("N70").Value = -101045
$ws.Range("H73").Value = 100459
$ws.Range("J73").Value = 100459
$ws.Range("L73").Value = 100459
$ws.Range("N73").Value = -102487
$ws.Range("H105").Value = 5717.8945
$ws.Range("I105").Value = 5853
$ws.Range("J105").Value = 5567.778
$ws.Range("K105").Value = 5853
$ws.Range("L105").Value = 5567.778
$ws.Range("M105").Value = -4106
$ws.Range("N105").Value = -9061.778

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 559
$ws.Range("I8").Value = 559
$ws.Range("K8").Value = 559
$ws.Range("M8").Value = -419
$ws.Range("H58").Value = 3127.8
$ws.Range("I58").Value = 796.1667
$ws.Range("J58").Value = 9123.429
$ws.Range("K58").Value = 796.1667
$ws.Range("L58").Value = 9123.429
$ws.Range("M58").Value = -593.1667
$ws.Range("N58").Value = -9529.429
$ws.Range("H62").Value = 7780
$ws.Range("I62").Value = 8925
$ws.Range("J62").Value = 5490
$ws.Range("K62").Value = 8925
$ws.Range("L62").Value = 5490
$ws.Range("M62").Value = -8301
$ws.Range("N62").Value = -6738
$ws.Range("H65").Value = 7780
$ws.Range("I65").Value = 8925
$ws.Range("J65").Value = 5490
$ws.Range("K65").Value = 44625
$ws.Range("L65").Value = 27450
$ws.Range("M65").Value = -41505
$ws.Range("N65").Value = -33690
$ws.Range("H132").Value = 1113.75
$ws.Range("I132").Value = 953.4808
$ws.Range("J132").Value = 3197.25
$ws.Range("K132").Value = 2860.4424
$ws.Range("L132").Value = 9591.75
$ws.Range("M132").Value = -330.4423999999999
$ws.Range("N132").Value = -14651.75
$ws.Range("H134").Value = 1497.2307
$ws.Range("I134").Value = 1335.8889
$ws.Range("J134").Value = 1860.25
$ws.Range("K134").Value = 4007.6667
$ws.Range("L134").Value = 5580.75
$ws.Range("M134").Value = -1472.6667
$ws.Range("N134").Value = -10650.75
$ws.Range("H136").Value = 3127.8
$ws.Range("I136").Value = 796.1667
$ws.Range("J136").Value = 9123.429
$ws.Range("K136").Value = 2388.5001
$ws.Range("L136").Value = 27370.287
$ws.Range("M136").Value = 161.4998999999998
$ws.Range("N136").Value = -32470.287

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1053
$ws.Range("I5").Value = 742.6667
$ws.Range("J5").Value = 1208.1666
$ws.Range("K5").Value = 2228.0001
$ws.Range("L5").Value = 3624.4998
$ws.Range("M5").Value = -2116.0001
$ws.Range("N5").Value = -3848.4998
$ws.Range("H97").Value = 466.66666
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -704
$ws.Range("N97").Value = -2792
$ws.Range("H113").Value = 678.6667
$ws.Range("J113").Value = 568.44446
$ws.Range("L113").Value = 1705.33338
$ws.Range("N113").Value = -6045.33338
$ws.Range("H122").Value = 715083.8
$ws.Range("J122").Value = 1000943.3
$ws.Range("L122").Value = 9008489.700000001
$ws.Range("N122").Value = -9013389.700000001
$ws.Range("H135").Value = 1053
$ws.Range("I135").Value = 742.6667
$ws.Range("J135").Value = 1208.1666
$ws.Range("K135").Value = 6684.0003
$ws.Range("L135").Value = 10873.4994
$ws.Range("M135").Value = -4149.0003
$ws.Range("N135").Value = -15943.4994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5001045
$ws.Range("I3").Value = 10001500
$ws.Range("J3").Value = 590
$ws.Range("K3").Value = 10001500
$ws.Range("L3").Value = 590
$ws.Range("M3").Value = -10001384
$ws.Range("N3").Value = -822
$ws.Range("H10").Value = 50000000
$ws.Range("I10").Value = 50000000
$ws.Range("K10").Value = 50000000
$ws.Range("M10").Value = -49999831
$ws.Range("H70").Value = 4018.524
$ws.Range("I70").Value = 3838.6667
$ws.Range("K70").Value = 3838.6667
$ws.Range("M70").Value = -3568.6667
$ws.Range("H73").Value = 4018.524
$ws.Range("I73").Value = 3838.6667
$ws.Range("K73").Value = 3838.6667
$ws.Range("M73").Value = -2902.6667
$ws.Range("H132").Value = 3376.1333
$ws.Range("I132").Value = 3414.2
$ws.Range("K132").Value = 10242.6
$ws.Range("M132").Value = -7712.599999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6737.048
$ws.Range("I122").Value = 9575.23
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 28725.69
$ws.Range("L122").Value = 6375
$ws.Range("M122").Value = -26275.69
$ws.Range("N122").Value = -11275

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2980
$ws.Range("J18").Value = 2980
$ws.Range("L18").Value = 2980
$ws.Range("N18").Value = -3326
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H62").Value = 6980.5
$ws.Range("I62").Value = 5400
$ws.Range("J62").Value = 8561
$ws.Range("K62").Value = 5400
$ws.Range("L62").Value = 8561
$ws.Range("M62").Value = -4776
$ws.Range("N62").Value = -9809
$ws.Range("H65").Value = 6980.5
$ws.Range("I65").Value = 5400
$ws.Range("J65").Value = 8561
$ws.Range("K65").Value = 27000
$ws.Range("L65").Value = 42805
$ws.Range("M65").Value = -23880
$ws.Range("N65").Value = -49045
$ws.Range("H122").Value = 1033.9333
$ws.Range("I122").Value = 1000.5714
$ws.Range("J122").Value = 1063.125
$ws.Range("K122").Value = 3001.7142
$ws.Range("L122").Value = 3189.375
$ws.Range("M122").Value = -551.7142000000003
$ws.Range("N122").Value = -8089.375
